$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; existing rows 20-116 shift down to 21-117.
$ws.Rows(20).Insert()

# Populate the newly inserted row 20 with the new daily price record.
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = 44749
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100107
$ws.Range("H20").Value = "Otros"
$ws.Range("I20").Value = 100107002
$ws.Range("J20").Value = "Chirimoya"
$ws.Range("K20").Value = "Cultivar IV Región"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = 5000
$ws.Range("O20").Value = 5000
$ws.Range("P20").Value = 5000
$ws.Range("Q20").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R20").Value = "Provincia del Elquí"
$ws.Range("S20").Value = 5000
$ws.Range("T20").Value = 1
